$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Paragraph 29 "Sendgrid: Mail API (Password via PM)" -> wrap "Sendgrid"
#    run with proofErr spellStart/spellEnd markers.
# ---------------------------------------------------------------------------
$xmlSendgrid = @"
<w:document $wns><w:body>
<w:p>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Sendgrid</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>: Mail API (Password via PM)</w:t></w:r>
</w:p>
</w:body></w:document>
"@
$null = $d.Paragraphs(29).Range.InsertXML($xmlSendgrid)

# ---------------------------------------------------------------------------
# 2) Insert a brand-new paragraph "Server Side eingeben Überprüfung (...)"
#    right after "Die Automatische Ausweitung..." (paragraph 15), before the
#    block of empty paragraphs.
# ---------------------------------------------------------------------------
$xmlAutomatische = @"
<w:document $wns><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r><w:t>Die Automatische Ausweitung von der über Subjekten implementieren</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Server Side eingeben Überprüfung (sicher </w:t></w:r>
  <w:r><w:t>gehen,</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t>dass</w:t></w:r>
  <w:r><w:t xml:space="preserve"> alle Client Informationen sin machen)</w:t></w:r>
</w:p>
</w:body></w:document>
"@
$null = $d.Paragraphs(15).Range.InsertXML($xmlAutomatische)

# ---------------------------------------------------------------------------
# 3) Paragraph 13 "E-Mail und Code überprüfen" -> add color rPr, and insert
#    two new paragraphs right after it: "Anmeldung durch E-Mail erlauben"
#    and "Code ablaufen lassen" (both ilvl 1, color 385623).
# ---------------------------------------------------------------------------
$xmlEMailCode = @"
<w:document $wns><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
    <w:t>E-Mail</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
    <w:t xml:space="preserve"> und Code überprüfen</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
    <w:t>Anmeldung durch E-Mail erlauben</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
    <w:t>Code ablaufen lassen</w:t>
  </w:r>
</w:p>
</w:body></w:document>
"@
$null = $d.Paragraphs(13).Range.InsertXML($xmlEMailCode)

# ---------------------------------------------------------------------------
# 4) Paragraph 11 "Registrierung Dialog hinzufügen" -> add color rPr (385623).
# ---------------------------------------------------------------------------
$xmlRegistrierung = @"
<w:document $wns><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
    <w:t>Registrierung Dialog hinzufügen</w:t>
  </w:r>
</w:p>
</w:body></w:document>
"@
$null = $d.Paragraphs(11).Range.InsertXML($xmlRegistrierung)

# ---------------------------------------------------------------------------
# 5) Paragraph 10 "Anmeldung nur durch Benedict E-Mail (ev. E-Mail-Bestätigung)"
#    -> add color rPr (70AD47) to the paragraph and each of its 3 runs.
# ---------------------------------------------------------------------------
$xmlAnmeldung = @"
<w:document $wns><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr>
    <w:t xml:space="preserve">Anmeldung nur durch Benedict E-Mail (ev. </w:t>
  </w:r>
  <w:r>
    <w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr>
    <w:t>E-Mail-Bestätigung</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr>
    <w:t>)</w:t>
  </w:r>
</w:p>
</w:body></w:document>
"@
$null = $d.Paragraphs(10).Range.InsertXML($xmlAnmeldung)

# ---------------------------------------------------------------------------
# 6) Paragraph 8 "Die Verschlüsselung an allen Queries s etablieren" -> split
#    "Queries" into its own spell-checked run.
# ---------------------------------------------------------------------------
$xmlQueries2 = @"
<w:document $wns><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
    <w:t xml:space="preserve">Die Verschlüsselung an allen </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
    <w:t>Queries</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
    <w:t xml:space="preserve"> s etablieren</w:t>
  </w:r>
</w:p>
</w:body></w:document>
"@
$null = $d.Paragraphs(8).Range.InsertXML($xmlQueries2)

# ---------------------------------------------------------------------------
# 7) Paragraph 5 "Die nötigen Informationen für die Entschlüsselung and den DB
#    Queries bekommen" -> split "Queries" into its own spell-checked run.
# ---------------------------------------------------------------------------
$xmlQueries1 = @"
<w:document $wns><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
    <w:t xml:space="preserve">Die nötigen Informationen für die Entschlüsselung and den DB </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
    <w:t>Queries</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>
    <w:t xml:space="preserve"> bekommen</w:t>
  </w:r>
</w:p>
</w:body></w:document>
"@
$null = $d.Paragraphs(5).Range.InsertXML($xmlQueries1)

# ---------------------------------------------------------------------------
# 8) Paragraph 3 "Token zu sha String Umwandeln" -> split "sha" into its own
#    spell-checked run.
# ---------------------------------------------------------------------------
$xmlToken = @"
<w:document $wns><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr>
    <w:t xml:space="preserve">Token zu </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr>
    <w:t>sha</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr>
    <w:t xml:space="preserve"> String Umwandeln</w:t>
  </w:r>
</w:p>
</w:body></w:document>
"@
$null = $d.Paragraphs(3).Range.InsertXML($xmlToken)

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
